# Update "想去人数" (F column) counts across sheets per the commit diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5334
$ws1.Range("F3").Value = 579
$ws1.Range("F4").Value = 11089
$ws1.Range("F5").Value = 276
$ws1.Range("F6").Value = 583
$ws1.Range("F7").Value = 157
$ws1.Range("F8").Value = 224
$ws1.Range("F9").Value = 947
$ws1.Range("F10").Value = 93

# Sheet "演出" (sheetId 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 23

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5334
$ws4.Range("F5").Value = 579
$ws4.Range("F6").Value = 23
$ws4.Range("F7").Value = 11089
$ws4.Range("F8").Value = 276
$ws4.Range("F9").Value = 583
$ws4.Range("F10").Value = 157
$ws4.Range("F13").Value = 224
$ws4.Range("F14").Value = 947
$ws4.Range("F16").Value = 93
